$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.856.09"
$ws.Range("E2").Value = "  +1.45%  "

$ws.Range("D3").Value = "1.768.66"
$ws.Range("E3").Value = "  +1.78%  "

$ws.Range("D4").Value = "'1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.37%  "

$ws.Range("D5").Value = "'328.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.07%  "

$ws.Range("E6").Value = "  -0.33%  "

$ws.Range("D7").Value = "'0.4468"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.77%  "

$ws.Range("D8").Value = "'0.3562"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.83%  "

$ws.Range("D9").Value = "'0.07444"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.73%  "

$ws.Range("D10").Value = "'42.03"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.05%  "

$ws.Range("D11").Value = "'1.095"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.09%  "

$ws.Range("D12").Value = "'1.001"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.38%  "

$ws.Range("D13").Value = "'20.94"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.53%  "

$ws.Range("D14").Value = "'6.023"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.71%  "

$ws.Range("D15").Value = "'7.235"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.45%  "

$ws.Range("D16").Value = "1.771.18"
$ws.Range("E16").Value = "  +1.61%  "

$ws.Range("D17").Value = "'93.38"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.41%  "

$ws.Range("D18").Value = "'0.00001060"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.51%  "

$ws.Range("D19").Value = "'0.06429"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.52%  "

$ws.Range("E20").Value = "  -0.20%  "

$ws.Range("D21").Value = "'17.10"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.75%  "

$ws.Range("E22").Value = "  +0.91%  "

$ws.Range("D23").Value = "27.908.25"
$ws.Range("E23").Value = "  +1.37%  "

$ws.Range("D24").Value = "'11.29"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.67%  "

$ws.Range("D25").Value = "'2.118"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.61%  "

$ws.Range("D26").Value = "'162.91"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.04%  "

$ws.Range("D27").Value = "'20.36"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.90%  "

$ws.Range("D28").Value = "1.974.33"
$ws.Range("E28").Value = "  +1.62%  "

$ws.Range("D29").Value = "'2.163"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.02%  "

$ws.Range("D30").Value = "'125.16"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.07%  "

$ws.Range("D31").Value = "'1.107"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.11%  "

$ws.Range("D32").Value = "'0.09179"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.05%  "

$ws.Range("D33").Value = "'5.599"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.30%  "

$ws.Range("D34").Value = "'3.651"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.46%  "

$ws.Range("D35").Value = "'11.86"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.29%  "

$ws.Range("D36").Value = "'0.02290"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.96%  "

$ws.Range("D37").Value = "'0.06118"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.02%  "

$ws.Range("E38").Value = "  +1.72%  "

$ws.Range("D39").Value = "'0.6323"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.21%  "

$ws.Range("D40").Value = "'4.959"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.06%  "

$ws.Range("D41").Value = "'1.188"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.42%  "

$ws.Range("D42").Value = "'1.393"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.15%  "

$ws.Range("D43").Value = "'7.910"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.86%  "

$ws.Range("D44").Value = "'13.22"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.72%  "

$ws.Range("D45").Value = "'3.741"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.18%  "

$ws.Range("D46").Value = "'0.5897"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.47%  "

$ws.Range("D47").Value = "'122.38"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.49%  "

$ws.Range("E48").Value = "  +1.62%  "

$ws.Range("E49").Value = "  +0.98%  "

$ws.Range("E50").Value = "  +1.35%  "

$ws.Range("D51").Value = "'72.90"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.45%  "
